# Update the "Dépenses" sheet row 2 with new values (unit testing data refresh
# on the legacy expense-tracking report).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "1971-08-27"

$ws.Range("A2").Value = "shopping with Degas"

$ws.Range("E2").Value = "Lightweight Wool Bench"
$ws.Range("F2").Value = "Tea"

$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "4.0"

$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "2.0"

# Writing numeric-looking text via .Value stamps a "Text" number format on
# the cell (new style). Restore the original (unformatted) style on the
# affected cells by pasting the formatting from an untouched neighbor that
# still carries the original row style.
$ws.Range("C2").Copy()
$ws.Range("B2").PasteSpecial(-4122)
$ws.Range("G2").PasteSpecial(-4122)
$ws.Range("H2").PasteSpecial(-4122)
$excel.CutCopyMode = $false
